# Update "想去人数" (wish-to-attend count) figures for several events across
# the workbook's sheets, matching the refreshed data pull recorded in the
# commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 8179   # was 8174
$wsExpo.Range("F5").Value  = 5959   # was 5955
$wsExpo.Range("F6").Value  = 506    # was 505
$wsExpo.Range("F7").Value  = 97     # was 96
$wsExpo.Range("F10").Value = 302    # was 301
$wsExpo.Range("F11").Value = 737    # was 714
$wsExpo.Range("F12").Value = 74     # was 72

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 92      # was 91

# --- Sheet "全部类型" (All types, combines the above) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 8179    # was 8174
$wsAll.Range("F5").Value  = 5959    # was 5955
$wsAll.Range("F6").Value  = 506     # was 505
$wsAll.Range("F7").Value  = 97      # was 96
$wsAll.Range("F10").Value = 302     # was 301
$wsAll.Range("F11").Value = 92      # was 91
$wsAll.Range("F15").Value = 737     # was 714
$wsAll.Range("F16").Value = 74      # was 72
